$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 235
$ws.Range("F4").Value = 824
$ws.Range("F6").Value = 406
$ws.Range("F7").Value = 569
$ws.Range("F8").Value = 216
$ws.Range("F11").Value = 136
$ws.Range("F12").Value = 618
$ws.Range("F14").Value = 1770
$ws.Range("F15").Value = 326
$ws.Range("F16").Value = 2615
$ws.Range("F17").Value = 303
$ws.Range("F18").Value = 488
$ws.Range("F19").Value = 46
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 222
$ws.Range("F7").Value = 473
$ws.Range("F13").Value = 86
$ws.Range("F14").Value = 38
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5298
$ws.Range("F3").Value = 309
$ws.Range("F4").Value = 211
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5298
$ws.Range("F4").Value = 309
$ws.Range("F6").Value = 211
$ws.Range("F7").Value = 235
$ws.Range("F8").Value = 222
$ws.Range("F12").Value = 473
$ws.Range("F13").Value = 824
$ws.Range("F17").Value = 406
$ws.Range("F18").Value = 569
$ws.Range("F19").Value = 216
$ws.Range("F20").Value = 59
$ws.Range("F23").Value = 136
$ws.Range("F26").Value = 618
$ws.Range("F28").Value = 86
$ws.Range("F29").Value = 1770
$ws.Range("F30").Value = 326
$ws.Range("F31").Value = 2615
$ws.Range("F32").Value = 38
$ws.Range("F33").Value = 303
$ws.Range("F34").Value = 488
$ws.Range("F35").Value = 46
